$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.745.79'
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").Value = '2.572.47'
$ws.Range("E3").Value = '  +0.53%  '

$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'" + '581.11'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").Value = "'" + '143.62'
$ws.Range("E6").Value = '  -2.95%  '

$ws.Range("D7").Value = "'" + '1.00'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +0.41%  '

$ws.Range("E9").Value = '  -2.47%  '

$ws.Range("D10").Value = "'" + '5.56'
$ws.Range("E10").Value = '  -1.05%  '

$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("E12").Value = '  -1.78%  '

$ws.Range("D13").Value = "'" + '26.81'
$ws.Range("E13").Value = '  -2.94%  '

$ws.Range("D14").Value = '3.032.33'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").Value = '62.695.96'
$ws.Range("E15").Value = '  -0.56%  '

$ws.Range("E16").Value = '  -2.18%  '

$ws.Range("D17").Value = '2.573.13'
$ws.Range("E17").Value = '  +0.67%  '

$ws.Range("D18").Value = "'" + '11.06'
$ws.Range("E18").Value = '  -2.72%  '

$ws.Range("D19").Value = "'" + '339.66'
$ws.Range("E19").Value = '  -0.72%  '

$ws.Range("E20").Value = '  -2.26%  '

$ws.Range("E21").Value = '  -2.60%  '

$ws.Range("D22").Value = "'" + '0.999'
$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("D23").Value = "'" + '66.89'
$ws.Range("E23").Value = '  +0.77%  '

$ws.Range("D24").Value = "'" + '1.57'
$ws.Range("E24").Value = '  -4.52%  '

$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = "'" + '0.164'
$ws.Range("E25").Value = '  -3.98%  '

$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").Value = "'" + '1.51'
$ws.Range("E26").Value = '  +1.88%  '

$ws.Range("D27").Value = "'" + '0.998'
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("D28").Value = "'" + '7.90'
$ws.Range("E28").Value = '  -2.64%  '

$ws.Range("E29").Value = '  -4.02%  '

$ws.Range("E30").Value = '  -3.39%  '

$ws.Range("D31").Value = "'" + '454.26'
$ws.Range("E31").Value = '  +3.17%  '

$ws.Range("E32").Value = '  -3.93%  '

$ws.Range("E33").Value = '  +0.99%  '

$ws.Range("D34").Value = "'" + '176.51'
$ws.Range("E34").Value = '  -0.27%  '

$ws.Range("D36").Value = "'" + '0.397'
$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("E38").Value = '  -1.77%  '

$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("E40").Value = '  -3.54%  '

$ws.Range("D41").Value = "'" + '40.12'
$ws.Range("E41").Value = '  +1.29%  '

$ws.Range("D42").Value = "'" + '157.15'
$ws.Range("E42").Value = '  +4.20%  '

$ws.Range("E43").Value = '  -3.94%  '

$ws.Range("D44").Value = "'" + '0.630'
$ws.Range("E44").Value = '  +2.88%  '

$ws.Range("D45").Value = "'" + '20.99'
$ws.Range("E45").Value = '  -0.40%  '

$ws.Range("E46").Value = '  -3.19%  '

$ws.Range("E47").Value = '  -2.02%  '

$ws.Range("E48").Value = '  -3.08%  '

$ws.Range("E49").Value = '  -2.75%  '

$ws.Range("E50").Value = '  +0.28%  '

$ws.Range("E51").Value = '  -4.12%  '

